# Hotfix: Fri Nov 15 16:22:21 RTZ 2024
#
# - Python sheet: add a "D" (updated-at) timestamp column (and a touched,
#   still-empty "E" column) to every existing row, and append a new row
#   documenting how to print the current date/time.
# - HTML sheet: the single placeholder/test row is removed.
# - CSS sheet: the placeholder/test row is replaced with real content and a
#   new row about the border-radius property is appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Python sheet
# ---------------------------------------------------------------------
$wsPython = $wb.Worksheets.Item("Python")

# Stamp every existing data row (1-35) with the same "updated at" value in
# column D. Column E stays logically empty, but we touch its number format
# so the cell (and therefore the sheet dimension) keeps being tracked, just
# like it is in the source workbook.
$wsPython.Range("D1:D35").Value = "2024-11-15 13:09:33"
$wsPython.Range("E1:E35").NumberFormat = "General"

# Append the new row (36) documenting printing the current date/time.
$pyCode = @'
import datetime
dt_now = datetime.datetime.now()
print(dt_now)
'@

$wsPython.Range("A36").Value = 2116
$wsPython.Range("B36").Value = $pyCode
$wsPython.Range("C36").Value = "Вывод текущей латы и времени"
$wsPython.Range("D36").Value = "2024-11-15 13:09:33"
$wsPython.Range("E36").NumberFormat = "General"

# Multi-line content auto-expands the row height; put it back to the
# worksheet's normal (non-custom) row height like the other rows.
$wsPython.Rows.Item(36).AutoFit()

# ---------------------------------------------------------------------
# HTML sheet
# ---------------------------------------------------------------------
$wsHtml = $wb.Worksheets.Item("HTML")
$wsHtml.Range("A1:E1").ClearContents()

# ---------------------------------------------------------------------
# CSS sheet
# ---------------------------------------------------------------------
$wsCss = $wb.Worksheets.Item("CSS")

# Replace the placeholder row with the real first CSS entry.
$wsCss.Range("A1").Value = 3
$wsCss.Range("B1").Value = "Первая запись в CSS11"
$wsCss.Range("C1").Value = "Первая запись в CSS11"
$wsCss.Range("D1").Value = "2024-11-15 10:41:09"
$wsCss.Range("E1").Value = "2024-11-15 14:16:23"

$cssBorderRadius = @'
/* Свойство задает округление элементу HTML */
border-radius: 10px;
/* top-left-and-bottom-right | top-right-and-bottom-left */
border-radius: 10px 5%;
/* top-left | top-right-and-bottom-left | bottom-right */
border-radius: 2px 4px 2px;
/* top-left | top-right | bottom-right | bottom-left */
border-radius: 1px 0 3px 4px;
'@

$wsCss.Range("A2").Value = 4
$wsCss.Range("B2").Value = "Свойство border-radius"
$wsCss.Range("C2").Value = $cssBorderRadius
$wsCss.Range("D2").Value = "2024-11-15 11:23:21"
$wsCss.Range("E2").NumberFormat = "General"

# Multi-line content auto-expands the row height; put it back to the
# worksheet's normal (non-custom) row height like row 1.
$wsCss.Rows.Item(2).AutoFit()
